$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Production (MW)" values for rows 2-31 (column B)
$newValues = @(646,611,564,524,467,436,421,404,370,344,340,349,375,408,411,399,391,399,397,378,358,363,384,389,390,393,391,389,368,354)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

# Shift every timestamp in column A (rows 2-97) forward by 15 days
for ($row = 2; $row -le 97; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $cell.Value2 + 15
}
